# Apply cryptos list price/volume refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.318.76"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.949.29"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.66"
$ws.Range("E5").Value = "  -2.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.10"
$ws.Range("E6").Value = "  +4.05%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").Value = "2.941.74"
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("E10").Value = "  -4.49%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.22"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D16").Value = "65.292.93"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "3.442.59"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.95"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "2.952.49"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.81"
$ws.Range("E20").Value = "  +8.01%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "446.19"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.25"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.10"
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("E27").Value = "  -5.90%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.03"
$ws.Range("E29").Value = "  +3.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.39"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.07"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.70"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.99"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "44.20"
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("E40").Value = "  -6.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.85"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.298"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.45"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "386.55"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0352"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "2.715.55"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.08"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.17"
$ws.Range("E50").Value = "  +5.55%  "
$ws.Range("E51").Value = "  +0.40%  "
